$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"3.324593"
$ws.Range("H2").Value = [double]"9.973779"
$ws.Range("I2").Value = [double]"0.06006369988653708"
$ws.Range("J2").Value = [double]"0.06006369988653707"
$ws.Range("M2").Value = [double]"3.303267"
$ws.Range("N2").Value = [double]"9.909801000000002"
$ws.Range("O2").Value = [double]"0.03362563178859915"
$ws.Range("P2").Value = [double]"0.03362563178859915"
$ws.Range("Q2").Value = [double]"10.982018345331"
$ws.Range("R2").Value = [double]"98.83816510797902"
$ws.Range("S2").Value = [double]"0.002019679856245621"
$ws.Range("T2").Value = [double]"0.00201967985624562"

# Row 3
$ws.Range("G3").Value = [double]"3.324593"
$ws.Range("H3").Value = [double]"9.973779"
$ws.Range("I3").Value = [double]"0.06006369988653708"
$ws.Range("J3").Value = [double]"0.06006369988653707"
$ws.Range("M3").Value = [double]"37.82684066666667"
$ws.Range("O3").Value = [double]"0.3850586149964086"
$ws.Range("P3").Value = [double]"0.3850586149964086"
$ws.Range("Q3").Value = [double]"125.7588496925153"
$ws.Range("R3").Value = [double]"1131.829647232638"
$ws.Range("S3").Value = [double]"0.02312804508986991"
$ws.Range("T3").Value = [double]"0.02312804508986991"

# Row 4
$ws.Range("G4").Value = [double]"3.324593"
$ws.Range("H4").Value = [double]"9.973779"
$ws.Range("I4").Value = [double]"0.06006369988653708"
$ws.Range("J4").Value = [double]"0.06006369988653707"
$ws.Range("M4").Value = [double]"9.149395999999999"
$ws.Range("N4").Value = [double]"27.448188"
$ws.Range("O4").Value = [double]"0.09313634682999644"
$ws.Range("P4").Value = [double]"0.09313634682999644"
$ws.Range("Q4").Value = [double]"30.418017895828"
$ws.Range("R4").Value = [double]"273.762161062452"
$ws.Range("S4").Value = [double]"0.005594113584525336"
$ws.Range("T4").Value = [double]"0.005594113584525335"

# Row 5
$ws.Range("G5").Value = [double]"3.324593"
$ws.Range("H5").Value = [double]"9.973779"
$ws.Range("I5").Value = [double]"0.06006369988653708"
$ws.Range("J5").Value = [double]"0.06006369988653707"
$ws.Range("M5").Value = [double]"47.95707433333333"
$ws.Range("N5").Value = [double]"143.871223"
$ws.Range("O5").Value = [double]"0.4881794063849957"
$ws.Range("P5").Value = [double]"0.4881794063849957"
$ws.Range("Q5").Value = [double]"159.4377536290797"
$ws.Range("R5").Value = [double]"1434.939782661717"
$ws.Range("S5").Value = [double]"0.0293218613558962"
$ws.Range("T5").Value = [double]"0.0293218613558962"

# Row 6
$ws.Range("I6").Value = [double]"0.3450317237488911"
$ws.Range("J6").Value = [double]"0.3450317237488911"
$ws.Range("M6").Value = [double]"3.303267"
$ws.Range("N6").Value = [double]"9.909801000000002"
$ws.Range("O6").Value = [double]"0.03362563178859915"
$ws.Range("P6").Value = [double]"0.03362563178859915"
$ws.Range("Q6").Value = [double]"63.085436413164"
$ws.Range("R6").Value = [double]"567.7689277184761"
$ws.Range("S6").Value = [double]"0.01160190969816587"
$ws.Range("T6").Value = [double]"0.01160190969816587"

# Row 7
$ws.Range("I7").Value = [double]"0.3450317237488911"
$ws.Range("J7").Value = [double]"0.3450317237488911"
$ws.Range("M7").Value = [double]"37.82684066666667"
$ws.Range("O7").Value = [double]"0.3850586149964086"
$ws.Range("P7").Value = [double]"0.3850586149964086"
$ws.Range("S7").Value = [double]"0.1328574376765715"
$ws.Range("T7").Value = [double]"0.1328574376765715"

# Row 8
$ws.Range("I8").Value = [double]"0.3450317237488911"
$ws.Range("J8").Value = [double]"0.3450317237488911"
$ws.Range("M8").Value = [double]"9.149395999999999"
$ws.Range("N8").Value = [double]"27.448188"
$ws.Range("O8").Value = [double]"0.09313634682999644"
$ws.Range("P8").Value = [double]"0.09313634682999644"
$ws.Range("Q8").Value = [double]"174.734176673232"
$ws.Range("R8").Value = [double]"1572.607590059088"
$ws.Range("S8").Value = [double]"0.03213499429042824"
$ws.Range("T8").Value = [double]"0.03213499429042824"

# Row 9
$ws.Range("I9").Value = [double]"0.3450317237488911"
$ws.Range("J9").Value = [double]"0.3450317237488911"
$ws.Range("M9").Value = [double]"47.95707433333333"
$ws.Range("N9").Value = [double]"143.871223"
$ws.Range("O9").Value = [double]"0.4881794063849957"
$ws.Range("P9").Value = [double]"0.4881794063849957"
$ws.Range("Q9").Value = [double]"915.8790262539719"
$ws.Range("R9").Value = [double]"8242.911236285747"
$ws.Range("S9").Value = [double]"0.1684373820837255"
$ws.Range("T9").Value = [double]"0.1684373820837255"

# Row 10
$ws.Range("G10").Value = [double]"32.86291466666667"
$ws.Range("H10").Value = [double]"98.58874400000001"
$ws.Range("I10").Value = [double]"0.5937172592060275"
$ws.Range("J10").Value = [double]"0.5937172592060274"
$ws.Range("M10").Value = [double]"3.303267"
$ws.Range("N10").Value = [double]"9.909801000000002"
$ws.Range("O10").Value = [double]"0.03362563178859915"
$ws.Range("P10").Value = [double]"0.03362563178859915"
$ws.Range("Q10").Value = [double]"108.554981542216"
$ws.Range("R10").Value = [double]"976.9948338799442"
$ws.Range("S10").Value = [double]"0.01996411794459816"
$ws.Range("T10").Value = [double]"0.01996411794459816"

# Row 11
$ws.Range("G11").Value = [double]"32.86291466666667"
$ws.Range("H11").Value = [double]"98.58874400000001"
$ws.Range("I11").Value = [double]"0.5937172592060275"
$ws.Range("J11").Value = [double]"0.5937172592060274"
$ws.Range("M11").Value = [double]"37.82684066666667"
$ws.Range("O11").Value = [double]"0.3850586149964086"
$ws.Range("P11").Value = [double]"0.3850586149964086"
$ws.Range("Q11").Value = [double]"1243.100236938263"
$ws.Range("R11").Value = [double]"11187.90213244437"
$ws.Range("S11").Value = [double]"0.2286159455293367"
$ws.Range("T11").Value = [double]"0.2286159455293366"

# Row 12
$ws.Range("G12").Value = [double]"32.86291466666667"
$ws.Range("H12").Value = [double]"98.58874400000001"
$ws.Range("I12").Value = [double]"0.5937172592060275"
$ws.Range("J12").Value = [double]"0.5937172592060274"
$ws.Range("M12").Value = [double]"9.149395999999999"
$ws.Range("N12").Value = [double]"27.448188"
$ws.Range("O12").Value = [double]"0.09313634682999644"
$ws.Range("P12").Value = [double]"0.09313634682999644"
$ws.Range("Q12").Value = [double]"300.6758199995413"
$ws.Range("R12").Value = [double]"2706.082379995872"
$ws.Range("S12").Value = [double]"0.05529665657236747"
$ws.Range("T12").Value = [double]"0.05529665657236746"

# Row 13
$ws.Range("G13").Value = [double]"32.86291466666667"
$ws.Range("H13").Value = [double]"98.58874400000001"
$ws.Range("I13").Value = [double]"0.5937172592060275"
$ws.Range("J13").Value = [double]"0.5937172592060274"
$ws.Range("M13").Value = [double]"47.95707433333333"
$ws.Range("N13").Value = [double]"143.871223"
$ws.Range("O13").Value = [double]"0.4881794063849957"
$ws.Range("P13").Value = [double]"0.4881794063849957"
$ws.Range("Q13").Value = [double]"1576.009241479323"
$ws.Range("R13").Value = [double]"14184.08317331391"
$ws.Range("S13").Value = [double]"0.2898405391597251"
$ws.Range("T13").Value = [double]"0.2898405391597251"

# Row 14
$ws.Range("E14").Value = [double]"2"
$ws.Range("F14").Value = [double]"0.6666666666666666"
$ws.Range("G14").Value = [double]"0.06571933333333334"
$ws.Range("H14").Value = [double]"0.197158"
$ws.Range("I14").Value = [double]"0.001187317158544407"
$ws.Range("J14").Value = [double]"0.001187317158544407"
$ws.Range("M14").Value = [double]"3.303267"
$ws.Range("N14").Value = [double]"9.909801000000002"
$ws.Range("O14").Value = [double]"0.03362563178859915"
$ws.Range("P14").Value = [double]"0.03362563178859915"
$ws.Range("Q14").Value = [double]"0.217088505062"
$ws.Range("R14").Value = [double]"1.953796545558"
$ws.Range("S14").Value = [double]"3.992428958950004E-05"
$ws.Range("T14").Value = [double]"3.992428958950003E-05"

# Row 15
$ws.Range("E15").Value = [double]"2"
$ws.Range("F15").Value = [double]"0.6666666666666666"
$ws.Range("G15").Value = [double]"0.06571933333333334"
$ws.Range("H15").Value = [double]"0.197158"
$ws.Range("I15").Value = [double]"0.001187317158544407"
$ws.Range("J15").Value = [double]"0.001187317158544407"
$ws.Range("M15").Value = [double]"37.82684066666667"
$ws.Range("O15").Value = [double]"0.3850586149964086"
$ws.Range("P15").Value = [double]"0.3850586149964086"
$ws.Range("Q15").Value = [double]"2.485954750719556"
$ws.Range("R15").Value = [double]"22.373592756476"
$ws.Range("S15").Value = [double]"0.0004571867006305807"
$ws.Range("T15").Value = [double]"0.0004571867006305806"

# Row 16
$ws.Range("E16").Value = [double]"2"
$ws.Range("F16").Value = [double]"0.6666666666666666"
$ws.Range("G16").Value = [double]"0.06571933333333334"
$ws.Range("H16").Value = [double]"0.197158"
$ws.Range("I16").Value = [double]"0.001187317158544407"
$ws.Range("J16").Value = [double]"0.001187317158544407"
$ws.Range("M16").Value = [double]"9.149395999999999"
$ws.Range("N16").Value = [double]"27.448188"
$ws.Range("O16").Value = [double]"0.09313634682999644"
$ws.Range("P16").Value = [double]"0.09313634682999644"
$ws.Range("Q16").Value = [double]"0.6012922055226667"
$ws.Range("R16").Value = [double]"5.411629849704"
$ws.Range("S16").Value = [double]"0.0001105823826753978"
$ws.Range("T16").Value = [double]"0.0001105823826753977"

# Row 17
$ws.Range("E17").Value = [double]"2"
$ws.Range("F17").Value = [double]"0.6666666666666666"
$ws.Range("G17").Value = [double]"0.06571933333333334"
$ws.Range("H17").Value = [double]"0.197158"
$ws.Range("I17").Value = [double]"0.001187317158544407"
$ws.Range("J17").Value = [double]"0.001187317158544407"
$ws.Range("M17").Value = [double]"47.95707433333333"
$ws.Range("N17").Value = [double]"143.871223"
$ws.Range("O17").Value = [double]"0.4881794063849957"
$ws.Range("P17").Value = [double]"0.4881794063849957"
$ws.Range("Q17").Value = [double]"3.151706953803778"
$ws.Range("R17").Value = [double]"28.365362584234"
$ws.Range("S17").Value = [double]"0.0005796237856489284"
$ws.Range("T17").Value = [double]"0.0005796237856489283"
